$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student"
$ws.Range("B1").Value = "Math"
$ws.Range("C1").Value = "Science"
$ws.Range("D1").Value = "English"
$ws.Range("E1").Value = "Gym"
